$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Investor *"
$ws.Range("B1").Value = "Fund *"
$ws.Range("C1").Value = "Committed Amount *"
$ws.Range("D1").Value = "Folio No *"

$ws.Range("D2").Select()
